{"js": "// Update the single filled-in data row of the commission table with the\n// new values described in the diff:\n//   06/05/2025                                   -> 06/05\n//   JACUTINGA - UTILAR DE JACUTINGA EIRELI EPP    -> Jacutinga - Utilar De Jacutinga Eireli Epp\n//   R$ 5.922,50                                   -> 5.922,50\n//   R$ -7,00                                      -> -7\n//   R$ 6.337,07                                   -> 5.507,93\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Find the data row: the row whose first cell currently holds the\n// \"06/05/2025\" date text (row 0 is the header, the remaining rows are\n// blank form-field template rows).\nrows.items.forEach((r) => r.load(\"values\"));\nawait context.sync();\n\nlet dataRowIndex = -1;\nfor (let i = 0; i < rows.items.length; i++) {\n  const vals = rows.items[i].values;\n  if (vals && vals[0] && vals[0][0] === \"06/05/2025\") {\n    dataRowIndex = i;\n    break;\n  }\n}\n\nif (dataRowIndex === -1) {\n  throw new Error(\"Could not locate the commission data row (06/05/2025).\");\n}\n\ntable.getCell(dataRowIndex, 0).value = \"06/05\";\ntable.getCell(dataRowIndex, 2).value = \"Jacutinga - Utilar De Jacutinga Eireli Epp\";\ntable.getCell(dataRowIndex, 4).value = \"5.922,50\";\ntable.getCell(dataRowIndex, 5).value = \"-7\";\ntable.getCell(dataRowIndex, 6).value = \"5.507,93\";\n\nawait context.sync();\n", "ps1": "# Update the single filled-in data row of the commission table with the\n# new values described in the diff:\n#   06/05/2025                                   -> 06/05\n#   JACUTINGA - UTILAR DE JACUTINGA EIRELI EPP    -> Jacutinga - Utilar De Jacutinga Eireli Epp\n#   R$ 5.922,50                                   -> 5.922,50\n#   R$ -7,00                                      -> -7\n#   R$ 6.337,07                                   -> 5.507,93\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Locate the data row: the row whose first cell currently holds the\n# \"06/05/2025\" date text (row 1 is the header, the remaining rows are\n# blank form-field template rows).\n$dataRow = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $txt = $t.Cell($r, 1).Range.Text\n    if ($txt -like \"06/05/2025*\") {\n        $dataRow = $r\n        break\n    }\n}\n\nif ($dataRow -eq 0) {\n    throw \"Could not locate the commission data row (06/05/2025).\"\n}\n\n$t.Cell($dataRow, 1).Range.Text = \"06/05\"\n$t.Cell($dataRow, 3).Range.Text = \"Jacutinga - Utilar De Jacutinga Eireli Epp\"\n$t.Cell($dataRow, 5).Range.Text = \"5.922,50\"\n$t.Cell($dataRow, 6).Range.Text = \"-7\"\n$t.Cell($dataRow, 7).Range.Text = \"5.507,93\"\n"}
